# This script applies the crypto price/volume/ranking updates described by the diff.
# Column D ("Price") values are numeric-looking strings (e.g. "315.13", "27.646.75") that
# must remain literal text (they use "." as a thousands grouping, not a decimal point, and
# some combine both groupings and decimals). A leading apostrophe forces Excel to store the
# value as text instead of auto-converting it to a number, matching the original inlineStr data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.646.75"
$ws.Range("E2").Value = "  -2.07%  "
$ws.Range("D3").Value = "'1.846.78"
$ws.Range("E3").Value = "  -0.93%  "
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("D5").Value = "'315.13"
$ws.Range("E5").Value = "  -1.29%  "
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("D7").Value = "'0.4253"
$ws.Range("E7").Value = "  -3.00%  "
$ws.Range("D8").Value = "'0.3656"
$ws.Range("E8").Value = "  -1.74%  "
$ws.Range("D9").Value = "'45.65"
$ws.Range("E9").Value = "  +1.12%  "
$ws.Range("D10").Value = "'0.07265"
$ws.Range("E10").Value = "  -3.61%  "
$ws.Range("E11").Value = "  -4.08%  "
$ws.Range("D12").Value = "'20.72"
$ws.Range("E12").Value = "  -2.77%  "
$ws.Range("D13").Value = "'1.823.53"
$ws.Range("E13").Value = "  -1.27%  "
$ws.Range("D14").Value = "'5.378"
$ws.Range("E14").Value = "  -1.30%  "
$ws.Range("D15").Value = "'6.568"
$ws.Range("E15").Value = "  -2.41%  "
$ws.Range("D16").Value = "'0.06837"
$ws.Range("E16").Value = "  -0.38%  "
$ws.Range("E17").Value = "  -0.07%  "
$ws.Range("D18").Value = "'78.11"
$ws.Range("E18").Value = "  -5.02%  "
$ws.Range("D19").Value = "'0.000008805"
$ws.Range("E19").Value = "  -3.43%  "
$ws.Range("E20").Value = "  -0.17%  "
$ws.Range("D21").Value = "'15.45"
$ws.Range("E21").Value = "  -3.50%  "
$ws.Range("D22").Value = "'27.628.13"
$ws.Range("E22").Value = "  -2.12%  "
$ws.Range("D23").Value = "'4.964"
$ws.Range("E23").Value = "  -3.88%  "
$ws.Range("D24").Value = "'10.61"
$ws.Range("E24").Value = "  -1.39%  "
$ws.Range("D25").Value = "'2.066.03"
$ws.Range("E25").Value = "  -1.32%  "
$ws.Range("E26").Value = "  +0.26%  "
$ws.Range("D27").Value = "'154.13"
$ws.Range("E27").Value = "  -0.47%  "
$ws.Range("D28").Value = "'18.27"
$ws.Range("E28").Value = "  -0.96%  "
$ws.Range("D29").Value = "'5.282"
$ws.Range("E29").Value = "  -1.48%  "
$ws.Range("D30").Value = "'1.816"
$ws.Range("E30").Value = "  +4.75%  "
$ws.Range("D31").Value = "'110.75"
$ws.Range("E31").Value = "  -3.14%  "
$ws.Range("D32").Value = "'0.08874"
$ws.Range("E32").Value = "  -2.05%  "
$ws.Range("D33").Value = "'0.7713"
$ws.Range("E33").Value = "  -3.98%  "
$ws.Range("D34").Value = "'4.560"
$ws.Range("E34").Value = "  -6.14%  "
$ws.Range("D35").Value = "'2.966"
$ws.Range("E35").Value = "  +0.44%  "
$ws.Range("E36").Value = "  -7.28%  "
$ws.Range("D37").Value = "'1.000"
$ws.Range("E37").Value = "  -0.30%  "
$ws.Range("D38").Value = "'0.05414"
$ws.Range("E38").Value = "  -0.95%  "
$ws.Range("D39").Value = "'1.097"
$ws.Range("E39").Value = "  -2.51%  "
$ws.Range("E40").Value = "  -1.09%  "
$ws.Range("D41").Value = "'2.897"
$ws.Range("E41").Value = "  -2.69%  "
$ws.Range("D42").Value = "'0.5071"
$ws.Range("E42").Value = "  -3.62%  "
$ws.Range("D43").Value = "'6.815"
$ws.Range("E43").Value = "  -4.65%  "
$ws.Range("D44").Value = "'0.1644"
$ws.Range("E44").Value = "  -1.79%  "
$ws.Range("D45").Value = "'8.238"
$ws.Range("E45").Value = "  -6.22%  "
$ws.Range("E46").Value = "  -2.02%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'10.38"
$ws.Range("E47").Value = "  -1.77%  "
$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D48").Value = "'0.4714"
$ws.Range("E48").Value = "  -3.38%  "
$ws.Range("D49").Value = "'105.54"
$ws.Range("E49").Value = "  -2.07%  "
$ws.Range("E50").Value = "  -0.23%  "
$ws.Range("D51").Value = "'1.641"
$ws.Range("E51").Value = "  -2.49%  "
